$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "75÷8=9, 3"
$t.Cell(1,2).Range.Text  = "13÷6=2, 1"
$t.Cell(1,3).Range.Text  = "75÷7=10, 5"
$t.Cell(1,4).Range.Text  = "27÷4=6, 3"
$t.Cell(1,5).Range.Text  = "85÷5=17, 0"

$t.Cell(5,1).Range.Text  = "90÷7=12, 6"
$t.Cell(5,2).Range.Text  = "43÷2=21, 1"
$t.Cell(5,3).Range.Text  = "72÷9=8, 0"
$t.Cell(5,5).Range.Text  = "66÷3=22, 0"

$t.Cell(9,1).Range.Text  = "82÷6=13, 4"
$t.Cell(9,2).Range.Text  = "23÷3=7, 2"
$t.Cell(9,3).Range.Text  = "54÷8=6, 6"
$t.Cell(9,4).Range.Text  = "33÷9=3, 6"
$t.Cell(9,5).Range.Text  = "87÷3=29, 0"

$t.Cell(13,1).Range.Text = "30÷7=4, 2"
$t.Cell(13,2).Range.Text = "83÷5=16, 3"
$t.Cell(13,3).Range.Text = "76÷9=8, 4"
$t.Cell(13,4).Range.Text = "91÷4=22, 3"
$t.Cell(13,5).Range.Text = "85÷8=10, 5"

$t.Cell(17,1).Range.Text = "77÷2=38, 1"
$t.Cell(17,2).Range.Text = "68÷8=8, 4"
$t.Cell(17,3).Range.Text = "61÷7=8, 5"
$t.Cell(17,4).Range.Text = "18÷6=3, 0"
$t.Cell(17,5).Range.Text = "55÷3=18, 1"
